$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.494.69'
$ws.Range("E2").Value = '  +1.52%  '
$ws.Range("D3").Value = '2.291.03'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '156.12'
$ws.Range("E5").Value = '  +15,494.57%  '
$ws.Range("D6").Value = '306.51'
$ws.Range("E6").Value = '  +0.07%  '
$ws.Range("D7").Value = '96.61'
$ws.Range("E7").Value = '  +4.22%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").Value = '0.497'
$ws.Range("E10").Value = '  +2.19%  '
$ws.Range("D11").Value = '35.56'
$ws.Range("E11").Value = '  +8.69%  '
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = '2.644.79'
$ws.Range("E15").Value = '  +0.64%  '
$ws.Range("D16").Value = '14.58'
$ws.Range("E16").Value = '  +1.80%  '
$ws.Range("D17").Value = '2.292.44'
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").Value = '0.798'
$ws.Range("E18").Value = '  +4.31%  '
$ws.Range("D19").Value = '42.341.53'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = '12.93'
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("D21").Value = '0.0₃0921'
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("D22").Value = '6.05'
$ws.Range("E22").Value = '  +1.68%  '
$ws.Range("D23").Value = '68.28'
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("D24").Value = '244.96'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").Value = '2.63'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '24.32'
$ws.Range("E28").Value = '  -0.08%  '
$ws.Range("D29").Value = '36.68'
$ws.Range("E29").Value = '  +6.77%  '
$ws.Range("D30").Value = '9.73'
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").Value = '2.11'
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("D32").Value = '161.37'
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("D33").Value = '5.39'
$ws.Range("E33").Value = '  +3.93%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '0.0756'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("E38").Value = '  +4.45%  '
$ws.Range("D39").Value = '2.39'
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").Value = '1.84'
$ws.Range("E40").Value = '  +1.04%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.116'
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("E42").Value = '  +7.65%  '
$ws.Range("D43").Value = '20.21'
$ws.Range("E43").Value = '  +3.85%  '
$ws.Range("D44").Value = '2.014.70'
$ws.Range("E44").Value = '  -2.87%  '
$ws.Range("D45").Value = '2.26'
$ws.Range("E45").Value = '  +10.58%  '
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("D47").Value = '10.29'
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("E48").Value = '  +2.88%  '
$ws.Range("D49").Value = '54.11'
$ws.Range("E49").Value = '  +4.24%  '
$ws.Range("E50").Value = '  +1.53%  '
$ws.Range("D51").Value = '73.18'
$ws.Range("E51").Value = '  +0.13%  '
